$p = $ppt.ActivePresentation
$p.Designs(1).Name = "Office Theme"
Write-Output $p.Designs(1).Name
